$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")

# Update the scenario cell B17 on ProductLoanInput: "RBI (India)" -> "Overdue/Due Fee/Int,Principal"
# (this also naturally drops the now-unused "RBI (India)" shared string and appends the
# new one, shifting every other shared-string index exactly like the target diff).
$wsInput.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Make ProductLoanInput the active sheet/tab (was ProductLoanOutput before the edit).
$wsInput.Activate()

# Move the selection on ProductLoanInput to B17 and scroll the view up so A7 is the
# top-left visible cell (was A28/B32 before the edit).
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
[void]$wsInput.Range("B17").Select()
